$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "iAU_TC_ID_30"
$ws.Range("B2").Value = "@RegressionA Validation of Questions list page(Negative Scenario) "
$ws.Range("C2").Value = "passed"
